$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.603.08"
$ws.Range("E2").Value = "  -3.32%  "
$ws.Range("D3").Value = "3.008.59"
$ws.Range("E3").Value = "  -2.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.008.25"
$ws.Range("E8").Value = "  -2.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E10").Value = "  -4.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.451"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "3.491.68"
$ws.Range("E15").Value = "  -3.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.111"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Value = "61.634.01"
$ws.Range("E17").Value = "  -3.11%  "
$ws.Range("D18").Value = "3.006.55"
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.681"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("E25").Value = "  -2.13%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("E34").Value = "  -5.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "458.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.38%  "
$ws.Range("D38").Value = "3.185.59"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0795"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("E40").Value = "  -3.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.119"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  -8.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.90%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.248"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.10%  "
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "120.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").Value = "0.0₃0501"
$ws.Range("E50").Value = "  -8.64%  "
$ws.Range("E51").Value = "  -1.16%  "
